# "Add files via upload" — the uploaded workbook adds four new bot-reply
# strings (about uploading/refusing to forward a chat transcript) to the
# tail of the "test" sheet (rows 75-78), which previously held empty
# placeholder cells. The sheet view's scroll position/selection also moved
# down to show the newly-added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")
$ws.Activate()

# Rows 75 & 77 hold long phrases -> wrap to two lines like the sheet's other
# multi-line entries (e.g. row 9), rows 76 & 78 are short one-liners.
$ws.Range("A75").Value = "Простите, вы пытаетесь отправить контактные данные, уберите их из сообщения и отправте его."
$ws.Range("A75").RowHeight = 30

$ws.Range("A76").Value = "Загрузить переписку"

$ws.Range("A77").Value = "Вот список переписок, выберите 1 из них и я пришлю вам текстовый фаил переписки:"
$ws.Range("A77").RowHeight = 30

$ws.Range("A78").Value = "Вот фаил переписки"

# The author's view had scrolled down and landed just past the new content.
$ws.Range("A80").Select() | Out-Null
